$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.6848073333333332
$ws.Range("H2").Value2 = 2.054422
$ws.Range("I2").Value2 = 0.2268310526442471
$ws.Range("J2").Value2 = 0.2268310526442472
$ws.Range("M2").Value2 = 19.58374133333333
$ws.Range("N2").Value2 = 58.751224
$ws.Range("O2").Value2 = 0.6578841466750758
$ws.Range("P2").Value2 = 0.6578841466750758
$ws.Range("Q2").Value2 = 13.41108967916978
$ws.Range("R2").Value2 = 120.699807112528
$ws.Range("S2").Value2 = 0.1492285535082697
$ws.Range("T2").Value2 = 0.1492285535082697
$ws.Range("G3").Value2 = 0.6848073333333332
$ws.Range("H3").Value2 = 2.054422
$ws.Range("I3").Value2 = 0.2268310526442471
$ws.Range("J3").Value2 = 0.2268310526442472
$ws.Range("M3").Value2 = 6.657374333333333
$ws.Range("O3").Value2 = 0.2236437337398222
$ws.Range("P3").Value2 = 0.2236437337398222
$ws.Range("Q3").Value2 = 4.559018764211777
$ws.Range("R3").Value2 = 41.031168877906
$ws.Range("S3").Value2 = 0.0507293435414936
$ws.Range("T3").Value2 = 0.05072934354149362
$ws.Range("G4").Value2 = 0.6848073333333332
$ws.Range("H4").Value2 = 2.054422
$ws.Range("I4").Value2 = 0.2268310526442471
$ws.Range("J4").Value2 = 0.2268310526442472
$ws.Range("M4").Value2 = 3.526650333333333
$ws.Range("N4").Value2 = 10.579951
$ws.Range("O4").Value2 = 0.118472119585102
$ws.Range("P4").Value2 = 0.118472119585102
$ws.Range("Q4").Value2 = 2.41507601036911
$ws.Range("R4").Value2 = 21.735684093322
$ws.Range("S4").Value2 = 0.02687315559448381
$ws.Range("T4").Value2 = 0.02687315559448382
$ws.Range("I5").Value2 = 0.1086184939966157
$ws.Range("J5").Value2 = 0.1086184939966157
$ws.Range("M5").Value2 = 19.58374133333333
$ws.Range("N5").Value2 = 58.751224
$ws.Range("O5").Value2 = 0.6578841466750758
$ws.Range("P5").Value2 = 0.6578841466750758
$ws.Range("Q5").Value2 = 6.421926569681777
$ws.Range("R5").Value2 = 57.797339127136
$ws.Range("S5").Value2 = 0.07145838523609535
$ws.Range("T5").Value2 = 0.07145838523609535
$ws.Range("I6").Value2 = 0.1086184939966157
$ws.Range("J6").Value2 = 0.1086184939966157
$ws.Range("M6").Value2 = 6.657374333333333
$ws.Range("O6").Value2 = 0.2236437337398222
$ws.Range("P6").Value2 = 0.2236437337398222
$ws.Range("R6").Value2 = 19.647855610972
$ws.Range("S6").Value2 = 0.0242918455505996
$ws.Range("T6").Value2 = 0.0242918455505996
$ws.Range("I7").Value2 = 0.1086184939966157
$ws.Range("J7").Value2 = 0.1086184939966157
$ws.Range("M7").Value2 = 3.526650333333333
$ws.Range("N7").Value2 = 10.579951
$ws.Range("O7").Value2 = 0.118472119585102
$ws.Range("P7").Value2 = 0.118472119585102
$ws.Range("Q7").Value2 = 1.156463879507111
$ws.Range("R7").Value2 = 10.408174915564
$ws.Range("S7").Value2 = 0.01286826320992073
$ws.Range("T7").Value2 = 0.01286826320992073
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 2.006290666666667
$ws.Range("H8").Value2 = 6.018872
$ws.Range("I8").Value2 = 0.6645504533591371
$ws.Range("J8").Value2 = 0.6645504533591372
$ws.Range("M8").Value2 = 19.58374133333333
$ws.Range("N8").Value2 = 58.751224
$ws.Range("O8").Value2 = 0.6578841466750758
$ws.Range("P8").Value2 = 0.6578841466750758
$ws.Range("Q8").Value2 = 39.29067745548089
$ws.Range("R8").Value2 = 353.616097099328
$ws.Range("S8").Value2 = 0.4371972079307107
$ws.Range("T8").Value2 = 0.4371972079307107
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 2.006290666666667
$ws.Range("H9").Value2 = 6.018872
$ws.Range("I9").Value2 = 0.6645504533591371
$ws.Range("J9").Value2 = 0.6645504533591372
$ws.Range("M9").Value2 = 6.657374333333333
$ws.Range("O9").Value2 = 0.2236437337398222
$ws.Range("P9").Value2 = 0.2236437337398222
$ws.Range("Q9").Value2 = 13.35662798947289
$ws.Range("R9").Value2 = 120.209651905256
$ws.Range("S9").Value2 = 0.148622544647729
$ws.Range("T9").Value2 = 0.148622544647729
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 2.006290666666667
$ws.Range("H10").Value2 = 6.018872
$ws.Range("I10").Value2 = 0.6645504533591371
$ws.Range("J10").Value2 = 0.6645504533591372
$ws.Range("M10").Value2 = 3.526650333333333
$ws.Range("N10").Value2 = 10.579951
$ws.Range("O10").Value2 = 0.118472119585102
$ws.Range("P10").Value2 = 0.118472119585102
$ws.Range("Q10").Value2 = 7.075485648363555
$ws.Range("R10").Value2 = 63.679370835272
$ws.Range("S10").Value2 = 0.07873070078069744
$ws.Range("T10").Value2 = 0.07873070078069745
